$d = $word.ActiveDocument

# The document currently has 3 paragraphs:
#   1) Title "Write Up"
#   2) empty paragraph  <-- this is the one we expand into 3 paragraphs
#   3) empty Heading1 paragraph (left untouched)
#
# Split paragraph 2 into three paragraphs by inserting two new paragraph
# marks after it, then fill in the text/style for each.

$target = $d.Paragraphs(2)
$target.Range.InsertParagraphAfter()
$target.Range.InsertParagraphAfter()

$d.Paragraphs(2).Range.Text = "This week, we will be looking at how we can add the wall to the room. The wall will be a tile, just like the floor is a tile, but we will be placing this wall on a separate layer from the floor, so there is a bit of a process that we will need to do in order to bring this wall tile in."

$d.Paragraphs(3).Range.Text = "So, if you would like to continue following along with this tutorial, and learn how to bring the wall in, then why don’t you join us for our brand-new article entitled:"

$d.Paragraphs(4).Range.Text = "5 The Wall"
$d.Paragraphs(4).Style = "Heading1"
